# "Yaw Angle Control Just Finished"
# Adds a new "Sheet2" after the existing "Sheet1", populates it with the
# yaw-angle control calculations, and updates the active-sheet/selection
# state on both sheets.

$wb = $excel.ActiveWorkbook

# Existing sheet (stays first).
$ws1 = $wb.Worksheets.Item(1)

# New sheet, inserted immediately after Sheet1 -> becomes sheets[1]/activeTab=1.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Yaw angle control data/formulas.
$ws2.Range("A1").Value = 3.352
$ws2.Range("A2").Value = 9.81
$ws2.Range("A3").Formula = "=A2*A1"

$ws2.Range("A5").Value = 25
$ws2.Range("B5").Formula = "=(A5-A3)/A1"

$ws2.Range("A6").Value = 50
$ws2.Range("B6").Formula = "=(A6-A3)/A1"

# Selection on Sheet1 moved from J12 -> J10 (and Sheet1 is no longer the
# tab-selected/active sheet).
$ws1.Range("J10").Select()

# Sheet2 becomes the active/tab-selected sheet with A5 selected.
$ws2.Range("A5").Select()
